$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.354.30"
$ws.Range("E2").Value = "  -4.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.568.30"
$ws.Range("E3").Value = "  -4.18%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.16"
$ws.Range("E6").Value = "  -3.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3685"
$ws.Range("E7").Value = "  -2.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.18"
$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3372"
$ws.Range("E9").Value = "  -4.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.163"
$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07595"
$ws.Range("E11").Value = "  -5.72%  "

$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.14"
$ws.Range("E13").Value = "  -3.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.046"
$ws.Range("E14").Value = "  -4.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.879"
$ws.Range("E15").Value = "  -5.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.570.31"
$ws.Range("E16").Value = "  -4.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001130"
$ws.Range("E17").Value = "  -5.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.92"
$ws.Range("E18").Value = "  -7.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06756"
$ws.Range("E19").Value = "  -2.79%  "

$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.229"
$ws.Range("E21").Value = "  -7.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5328"
$ws.Range("E22").Value = "  -6.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.49"
$ws.Range("E23").Value = "  -4.81%  "

$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.369.85"
$ws.Range("E25").Value = "  -4.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.383"
$ws.Range("E26").Value = "  -3.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.977"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.89"
$ws.Range("E28").Value = "  -4.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "144.96"
$ws.Range("E29").Value = "  -5.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.959"
$ws.Range("E30").Value = "  -4.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.14"
$ws.Range("E31").Value = "  -5.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.745.88"
$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.031"
$ws.Range("E33").Value = "  +5.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.238"
$ws.Range("E34").Value = "  -8.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.992"
$ws.Range("E35").Value = "  -6.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.24"
$ws.Range("E36").Value = "  -9.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08463"
$ws.Range("E37").Value = "  -2.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02520"
$ws.Range("E38").Value = "  -6.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2325"
$ws.Range("E39").Value = "  -3.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.520"
$ws.Range("E40").Value = "  -5.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06504"
$ws.Range("E41").Value = "  -4.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.79"
$ws.Range("E42").Value = "  -9.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.240"
$ws.Range("E43").Value = "  -4.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6349"
$ws.Range("E44").Value = "  -7.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.17"
$ws.Range("E45").Value = "  -8.16%  "

$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5969"
$ws.Range("E47").Value = "  -5.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.757"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.119"
$ws.Range("E49").Value = "  -5.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.264"
$ws.Range("E50").Value = "  +7.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.02"
$ws.Range("E51").Value = "  -3.02%  "
